$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E8").Select()
$ws.Range("E8").Value = "GIT UPDATE"
